$d = $word.ActiveDocument

$replacements = @(
    @("2024-10-07 Monday", "2024-10-08 Tuesday"),
    @("29÷2=14, 1", "88÷9=9, 7"),
    @("16÷5=3, 1", "29÷9=3, 2"),
    @("38÷9=4, 2", "96÷2=48, 0"),
    @("11÷7=1, 4", "82÷2=41, 0"),
    @("86÷7=12, 2", "91÷7=13, 0"),
    @("53÷4=13, 1", "21÷6=3, 3"),
    @("74÷4=18, 2", "83÷4=20, 3"),
    @("74÷3=24, 2", "76÷7=10, 6"),
    @("18÷3=6, 0", "99÷9=11, 0"),
    @("28÷6=4, 4", "53÷5=10, 3"),
    @("35÷4=8, 3", "18÷7=2, 4"),
    @("26÷7=3, 5", "25÷4=6, 1"),
    @("44÷9=4, 8", "25÷6=4, 1"),
    @("28÷3=9, 1", "37÷4=9, 1"),
    @("46÷2=23, 0", "55÷3=18, 1"),
    @("80÷2=40, 0", "17÷6=2, 5"),
    @("83÷8=10, 3", "26÷9=2, 8"),
    @("94÷5=18, 4", "40÷3=13, 1"),
    @("26÷8=3, 2", "78÷2=39, 0"),
    @("64÷3=21, 1", "78÷3=26, 0"),
    @("93÷4=23, 1", "37÷5=7, 2"),
    @("31÷3=10, 1", "98÷7=14, 0"),
    @("97÷3=32, 1", "59÷2=29, 1"),
    @("54÷4=13, 2", "78÷9=8, 6"),
    @("41÷8=5, 1", "20÷6=3, 2")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
